$wb = $excel.ActiveWorkbook

$updates = @{
    2  = 1729
    3  = 234
    4  = 222
    5  = 7171
    6  = 453
    7  = 486
    10 = 8993
    11 = 2374
    12 = 289
    13 = 9145
    14 = 10427
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
